# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the related "Latest HO Xliff Generate Date" / "Latest Handoff
# Datetime" timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns (zh-cn / de-de) and the handoff date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 22:49:42"

# zh-cn sheet: status + its own handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 22:49:37"

# de-de sheet: status + its own handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 22:49:42"

# The status text grew longer ("In Translation" -> "Ready for handoff"),
# so the columns that display it widen to fit the new content. The COM
# ColumnWidth property is quantized to 1/6-character steps by the host, so
# 16.333333333333332 is the closest reachable value to the recorded
# target stored width (~17.216 "characters").
$targetColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
